$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before D (which holds "Observacao"), shifting it to E
# and shifting "Responsavel" from E to F. This adds the new "Valor_Saida" field.
$ws.Columns("D").Insert()
$ws.Range("D1").Value = "Valor_Saida"

# Remove the old sample rows 3-7, keeping only the header and a single data row.
$ws.Rows("3:7").Delete()

# Update the remaining data row (row 2) with the new entry values.
$ws.Range("A2").Value = "23/02/2026 01:16:00"
$ws.Range("B2").Value = "Culto de Ceia"
$ws.Range("C2").Value = "2700,00"
$ws.Range("D2").Value = "100,00"
$ws.Range("E2").Value = ""
$ws.Range("F2").Value = "ADMINISTRADOR"
